# Swap the full data rows (columns B:AD, i.e. 2..30) between each pair of
# rows below. Column A (the running index) stays put; every other field
# (match id, teams, scores, odds, P/L columns, ...) moves with its row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(67, 68),
    @(173, 174),
    @(178, 179),
    @(188, 189),
    @(346, 347),
    @(384, 385),
    @(386, 387)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = 2; $col -le 30; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
